$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# These Price cells are numeric-looking text (e.g. "290.28") in the source
# data. Mark them as Text first so Excel keeps them as strings instead of
# auto-converting to numbers.
$ws.Range("D6","D7","D8","D9","D10","D11","D13","D15","D17","D18","D19","D21","D22","D23","D24","D26","D27","D28","D29","D31","D32","D33","D34","D35","D36","D37","D38","D39","D40","D41","D44","D45","D46","D47","D48","D49","D50","D51").NumberFormat = "@"

$ws.Range("D2").Value = '21.763.36'
$ws.Range("E2").Value = '  -1.66%  '
$ws.Range("D3").Value = '1.539.31'
$ws.Range("E3").Value = '  -1.47%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("E5").Value = '  +0.08%  '
$ws.Range("D6").Value = '290.28'
$ws.Range("E6").Value = '  +0.24%  '
$ws.Range("D7").Value = '0.3882'
$ws.Range("E7").Value = '  +2.02%  '
$ws.Range("D8").Value = '0.3190'
$ws.Range("E8").Value = '  -3.12%  '
$ws.Range("D9").Value = '43.17'
$ws.Range("E9").Value = '  -0.70%  '
$ws.Range("D10").Value = '0.07194'
$ws.Range("E10").Value = '  -2.52%  '
$ws.Range("D11").Value = '1.060'
$ws.Range("E11").Value = '  -7.12%  '
$ws.Range("E12").Value = '  +0.05%  '
$ws.Range("D13").Value = '5.632'
$ws.Range("E13").Value = '  -3.40%  '
$ws.Range("E14").Value = '  -7.27%  '
$ws.Range("D15").Value = '6.624'
$ws.Range("E15").Value = '  -4.04%  '
$ws.Range("D16").Value = '1.544.51'
$ws.Range("E16").Value = '  -1.68%  '
$ws.Range("D17").Value = '0.00001105'
$ws.Range("E17").Value = '  +0.87%  '
$ws.Range("D18").Value = '0.06571'
$ws.Range("E18").Value = '  -1.23%  '
$ws.Range("D19").Value = '83.43'
$ws.Range("E19").Value = '  -2.85%  '
$ws.Range("E20").Value = '  +0.09%  '
$ws.Range("D21").Value = '6.141'
$ws.Range("E21").Value = '  -4.99%  '
$ws.Range("D22").Value = '15.37'
$ws.Range("E22").Value = '  -4.82%  '
$ws.Range("D23").Value = '10.91'
$ws.Range("E23").Value = '  -6.86%  '
$ws.Range("D24").Value = '2.377'
$ws.Range("E24").Value = '  +4.64%  '
$ws.Range("D25").Value = '21.766.64'
$ws.Range("E25").Value = '  -1.71%  '
$ws.Range("D26").Value = '2.388'
$ws.Range("E26").Value = '  -6.65%  '
$ws.Range("D27").Value = '145.81'
$ws.Range("E27").Value = '  -3.45%  '
$ws.Range("D28").Value = '18.37'
$ws.Range("E28").Value = '  -4.01%  '
$ws.Range("D29").Value = '4.849'
$ws.Range("E29").Value = '  -0.48%  '
$ws.Range("D30").Value = '1.717.50'
$ws.Range("E30").Value = '  -1.12%  '
$ws.Range("D31").Value = '117.46'
$ws.Range("E31").Value = '  -3.35%  '
$ws.Range("D32").Value = '0.9681'
$ws.Range("E32").Value = '  -14.09%  '
$ws.Range("D33").Value = '5.884'
$ws.Range("E33").Value = '  -2.88%  '
$ws.Range("D34").Value = '0.08217'
$ws.Range("E34").Value = '  +0.35%  '
$ws.Range("D35").Value = '8.937'
$ws.Range("E35").Value = '  -4.77%  '
$ws.Range("D36").Value = '0.06090'
$ws.Range("E36").Value = '  -2.25%  '
$ws.Range("D37").Value = '5.129'
$ws.Range("E37").Value = '  -3.34%  '
$ws.Range("D38").Value = '1.483'
$ws.Range("E38").Value = '  -19.72%  '
$ws.Range("D39").Value = '0.02203'
$ws.Range("D40").Value = '0.2040'
$ws.Range("E40").Value = '  -4.88%  '
$ws.Range("D41").Value = '1.185'
$ws.Range("E41").Value = '  -4.47%  '
$ws.Range("E42").Value = '  +0.01%  '
$ws.Range("E43").Value = '  -3.82%  '
$ws.Range("D44").Value = '0.5748'
$ws.Range("E44").Value = '  -4.25%  '
$ws.Range("D45").Value = '13.09'
$ws.Range("E45").Value = '  -5.24%  '
$ws.Range("D46").Value = '3.747'
$ws.Range("E46").Value = '  -0.08%  '
$ws.Range("D47").Value = '0.5511'
$ws.Range("E47").Value = '  -5.05%  '
$ws.Range("D48").Value = '118.11'
$ws.Range("E48").Value = '  -2.25%  '
$ws.Range("D49").Value = '1.865'
$ws.Range("E49").Value = '  -5.87%  '
$ws.Range("D50").Value = '1.142'
$ws.Range("E50").Value = '  -2.65%  '
$ws.Range("D51").Value = '0.06734'
$ws.Range("E51").Value = '  -3.62%  '
